# Update cryptocurrency price/volume data in the worksheet to reflect
# the latest GitHub Actions scrape (commit: "Updated cryptos list on
# Thu Mar 14 13:23:28 UTC 2024 with GitHub Actions").
#
# Numeric-looking text values in column D must stay TEXT (the sheet
# stores prices as inline strings, e.g. "0.999" or "1.00"), so for those
# cells we temporarily force a text NumberFormat before assigning the
# value, then clear the format again so no stray cell style is left
# behind (matches original cells, which carry no explicit style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "72.592.73"
$ws.Range("E2").Value = "  -0.41%  "

# Row 3
$ws.Range("D3").Value = "3.932.05"
$ws.Range("E3").Value = "  -2.03%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.10"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.29%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.41"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +11.33%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.682"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.94%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.780"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.83%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.187"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +9.66%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "55.53"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.49%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000330"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.94%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.47"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +4.36%  "

# Row 14
$ws.Range("D14").Value = "4.552.82"
$ws.Range("E14").Value = "  -2.30%  "

# Row 15
$ws.Range("D15").Value = "3.931.60"
$ws.Range("E15").Value = "  -2.36%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.33"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.19%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.14"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.70%  "

# Row 18
$ws.Range("E18").Value = "  -3.98%  "

# Row 19
$ws.Range("D19").Value = "72.478.36"
$ws.Range("E19").Value = "  -0.37%  "

# Row 20
$ws.Range("E20").Value = "  -1.28%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "447.88"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.21%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.79"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.60%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "95.48"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.64%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.31"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.83%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.01"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.72%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.25"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.09%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.08"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.97%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.95"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.55%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.31"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.32%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.73"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.69%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.85"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.90%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.81"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.63%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.60"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.48%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.127"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.64%  "

# Row 35
$ws.Range("D35").Value = "0.0₃0984"
$ws.Range("E35").Value = "  +13.17%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "68.86"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.21%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "624.01"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -8.95%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.425"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.07%  "

# Row 39
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.16%  "

# Row 40
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.33"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.52%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.145"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.75%  "

# Row 42
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.07%  "

# Row 43
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.23"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +41.96%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0476"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.37%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.49"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.29%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.147"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.53%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.37"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.18%  "

# Row 48
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.58"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -7.23%  "

# Row 49
$ws.Range("E49").Value = "  -16.47%  "

# Row 50
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.826.97"
$ws.Range("E50").Value = "  -0.59%  "

# Row 51
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000276"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.01%  "
